$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 34 (shifting existing rows 34..129 down to 35..130)
$ws.Rows.Item(34).Insert()

# Fill in the constant columns (same for every record in this dataset)
$ws.Cells.Item(34, 1).Value = 7
$ws.Cells.Item(34, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(34, 3).Value = "Ñuble"
$ws.Cells.Item(34, 4).Value = 45099
$ws.Cells.Item(34, 5).Value = 16
$ws.Cells.Item(34, 6).Value = 100112030
$ws.Cells.Item(34, 7).Value = "Poroto granado"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 30
$ws.Cells.Item(34, 11).Value = 25000
$ws.Cells.Item(34, 12).Value = 25000
$ws.Cells.Item(34, 13).Value = 25000
$ws.Cells.Item(34, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(34, 16).Value = 1000
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"
